# Update "想去人数" (F column) values across the four sheets of the
# 广州-漫展信息 workbook, per the generated-output refresh commit.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> hashtable of row -> new F value
$updates = @{
    "展览" = @{
        2  = 2229
        3  = 254
        4  = 161
        5  = 155
        6  = 292
        8  = 665
        9  = 488
        10 = 594
        12 = 55
        14 = 941
        15 = 203
        16 = 124
        17 = 81
        20 = 201
        21 = 78
    }
    "演出" = @{
        2  = 62
        3  = 30
        6  = 162
        8  = 2404
        16 = 2223
    }
    "本地生活" = @{
        3 = 314
        4 = 155
    }
    "全部类型" = @{
        3  = 62
        4  = 30
        5  = 2229
        6  = 314
        7  = 254
        8  = 161
        9  = 155
        10 = 292
        14 = 162
        15 = 155
        16 = 665
        17 = 488
        18 = 594
        20 = 55
        22 = 941
        24 = 2404
        30 = 203
        31 = 124
        32 = 81
        37 = 201
        38 = 78
        39 = 2223
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
